$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format so numeric-looking strings
# (e.g. "46.687.04", "1.00") are preserved exactly as typed, not
# silently coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '46.687.04'
$ws.Range('E2').Value = '  +4.64%  '
$ws.Range('D3').Value = '2.333.82'
$ws.Range('E3').Value = '  +3.92%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.83%  '
$ws.Range('D5').Value = '307.14'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').Value = '98.34'
$ws.Range('E6').Value = '  +3.48%  '
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '0.539'
$ws.Range('E9').Value = '  +3.62%  '
$ws.Range('D10').Value = '36.30'
$ws.Range('E10').Value = '  +3.15%  '
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('E12').Value = '  +3.37%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '2.688.24'
$ws.Range('E14').Value = '  +3.94%  '
$ws.Range('D15').Value = '2.333.91'
$ws.Range('E15').Value = '  +4.30%  '
$ws.Range('E16').Value = '  +4.30%  '
$ws.Range('D17').Value = '0.834'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '46.527.31'
$ws.Range('E18').Value = '  +4.74%  '
$ws.Range('D19').Value = '12.96'
$ws.Range('E19').Value = '  +9.47%  '
$ws.Range('D20').Value = '0.0₃0950'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').Value = '6.20'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = '66.80'
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('D23').Value = '244.65'
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D27').Value = '41.97'
$ws.Range('E27').Value = '  +12.78%  '
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').Value = '9.85'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').Value = '20.17'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('D32').Value = '152.36'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').Value = '0.0808'
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('D34').Value = '2.63'
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').Value = '3.01'
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('D38').Value = '1.81'
$ws.Range('E38').Value = '  -2.99%  '
$ws.Range('D39').Value = '4.04'
$ws.Range('E39').Value = '  +6.60%  '
$ws.Range('E40').Value = '  +6.24%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').Value = '14.03'
$ws.Range('E42').Value = '  -8.30%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.97'
$ws.Range('E44').Value = '  +10.27%  '
$ws.Range('D45').Value = '1.792.79'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('E46').Value = '  +5.84%  '
$ws.Range('D47').Value = '74.68'
$ws.Range('E47').Value = '  +8.71%  '
$ws.Range('D48').Value = '81.25'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '55.70'
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '4.89'
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '98.19'
$ws.Range('E51').Value = '  -0.35%  '

# Restore the default (unstyled) cell style now that values are committed,
# matching the original workbook which has no explicit style on these cells.
$ws.Range("D2:E51").Style = "Normal"
